# Présentation Procédures stockées / Exercices Fil Rouge
#
# Add the new "Dépôts BDD FIL ROUGE: 16h" note to the Vendredi (Friday)
# column of the weekly schedule, and leave the selection on G11 (matching
# the author's cursor position when the file was saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = "Dépôts BDD FIL ROUGE: 16h"

$ws.Range("G11").Select()
